# Applies the changes described by the commit:
#  - "Agendar uma consulta" gains a suffix and absorbs/removes the
#    following "Verificar disponibilidade" paragraph.
#  - "chat bot" -> "chatbot" typo fix.
#  - Run-merge + typo fixes: "quem" -> "que", "previnir" -> "prevenir".
#  - Trim a long run-on sentence about the hospital feature.
#  - Restructure the "Licenças de fonte" section: add a blank paragraph,
#    blank out the "Licenças de fonte" paragraph's run, and delete the
#    "Licença específica da fonte..." paragraph entirely.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "Agendar uma consulta" -> "...consulta(Informações o hospital)" and
#    drop the whole next paragraph ("Verificar disponibilidade").
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Agendar uma consulta", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Agendar uma consulta(Informações o hospital)", 2)

$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    if ($paras.Item($i).Range.Text -eq "Verificar disponibilidade`r") {
        $paras.Item($i).Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------
# 2) "chat bot" -> "chatbot"
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "uma interação com chat bot que facilite", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "uma interação com chatbot que facilite", 2)

# ---------------------------------------------------------------------
# 3) Merge "s estudantes, pessoas prevenidas e doentes, " + "quem" +
#    " querem estudar/entender/previnir sobre determinada doença" into
#    a single corrected run.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "s estudantes, pessoas prevenidas e doentes, quem querem estudar/entender/previnir sobre determinada doença",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "s estudantes, pessoas prevenidas e doentes, que querem estudar/entender/prevenir sobre determinada doença",
    2)

# ---------------------------------------------------------------------
# 4) Trim the run-on sentence about the hospital functionality.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Ir na funcionalidade  hospital, procurar o hospital que deseja, e selecionar o hospital, com isso o app informa as informações do hospital, se ele quiser marcar a consulta direto no app ele pode, se não ele pega o caminho que a API oferecer do hospital",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Ir na funcionalidade  hospital, procurar o hospital que deseja, e selecionar o hospital, com isso o app informa as informações do hospital.",
    2)

# ---------------------------------------------------------------------
# 5) Restructure the "Licenças de fonte" section.
# ---------------------------------------------------------------------
$paras = $d.Paragraphs
$emGeralIdx = 0
$licencasIdx = 0
$licencaEspecificaIdx = 0
for ($i = 1; $i -le $paras.Count; $i++) {
    $t = $paras.Item($i).Range.Text
    if ($t -eq "Em geral o nosso app será usado previamente pelos clientes cadastrados, sem isso, o usuário não conseguirá usar o nosso app.`r") {
        $emGeralIdx = $i
    }
    if ($t -eq "Licenças de fonte`r") {
        $licencasIdx = $i
    }
    if ($t -eq "Licença específica da fonte. Relacionado às fontes usadas on-line ou internamente por uma organização.`r") {
        $licencaEspecificaIdx = $i
    }
}

# 5a) Insert a new, blank paragraph right after "Em geral...". We
#     replace the whole paragraph (original text + a fresh blank
#     paragraph) so the insertion point's inherited run formatting
#     does not leak into the new paragraph's run.
$emGeralPara = $paras.Item($emGeralIdx)
$emGeralXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:after="280" w:before="240" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Verdana" w:cs="Verdana" w:eastAsia="Verdana" w:hAnsi="Verdana"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:cs="Verdana" w:eastAsia="Verdana" w:hAnsi="Verdana"/><w:sz w:val="19"/><w:szCs w:val="19"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Em geral o nosso app será usado previamente pelos clientes cadastrados, sem isso, o usuário não conseguirá usar o nosso app.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="280" w:before="240" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Verdana" w:cs="Verdana" w:eastAsia="Verdana" w:hAnsi="Verdana"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p>'
$emGeralPara.Range.InsertXML($emGeralXml)

# 5b) Blank out the "Licenças de fonte" paragraph's run (keep its pPr
#     untouched, shrink the run to an empty, minimally-formatted one).
$paras = $d.Paragraphs
$licencasIdx = 0
for ($i = 1; $i -le $paras.Count; $i++) {
    if ($paras.Item($i).Range.Text -eq "Licenças de fonte`r") {
        $licencasIdx = $i
        break
    }
}
$licencasPara = $paras.Item($licencasIdx)
$licencasXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:keepNext w:val="0"/><w:keepLines w:val="0"/><w:widowControl w:val="1"/><w:pBdr><w:top w:space="0" w:sz="0" w:val="nil"/><w:left w:space="0" w:sz="0" w:val="nil"/><w:bottom w:space="0" w:sz="0" w:val="nil"/><w:right w:space="0" w:sz="0" w:val="nil"/><w:between w:space="0" w:sz="0" w:val="nil"/></w:pBdr><w:shd w:fill="auto" w:val="clear"/><w:spacing w:after="280" w:before="240" w:line="240" w:lineRule="auto"/><w:ind w:left="0" w:right="0" w:firstLine="0"/><w:jc w:val="left"/><w:rPr><w:rFonts w:ascii="Verdana" w:cs="Verdana" w:eastAsia="Verdana" w:hAnsi="Verdana"/><w:b w:val="1"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p>'
$licencasPara.Range.InsertXML($licencasXml)

# 5c) Delete the "Licença específica da fonte..." paragraph entirely.
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    if ($paras.Item($i).Range.Text -eq "Licença específica da fonte. Relacionado às fontes usadas on-line ou internamente por uma organização.`r") {
        $paras.Item($i).Range.Delete()
        break
    }
}

Write-Output "Edits applied."
